$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the existing row 621 (which holds 2026/12/29 火 13),
# shifting it and everything below down by one row (old 621..662 -> new 622..663).
$ws.Rows.Item(621).Insert()

# Populate the newly inserted row 621 with the new data point for 2026/01/12.
$ws.Cells.Item(621, 1).Value = "'2026/01/12"
$ws.Cells.Item(621, 1).Style = "Normal"
$ws.Cells.Item(621, 2).Value = "月"
$ws.Cells.Item(621, 3).Value = 16
$ws.Cells.Item(621, 4).Value = 201
